# Refresh cryptos list values (prices & 1h volume %) to match latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.495.49"
$ws.Range("E2").Value = "  +4.34%  "
$ws.Range("D3").Value = "4.038.67"
$ws.Range("E3").Value = "  +3.24%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "518.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.734"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +19.81%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.762"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.25%  "
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000327"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +11.65%  "
$ws.Range("E13").Value = "  +6.39%  "
$ws.Range("D14").Value = "4.688.19"
$ws.Range("E14").Value = "  +3.35%  "
$ws.Range("D15").Value = "4.054.89"
$ws.Range("E15").Value = "  +3.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("E19").Value = "  -1.56%  "
$ws.Range("D20").Value = "72.360.21"
$ws.Range("E20").Value = "  +4.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "448.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "104.87"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +18.58%  "
$ws.Range("E23").Value = "  +5.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.20%  "
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.34%  "
$ws.Range("E30").Value = "  +9.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.65"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.64%  "
$ws.Range("E32").Value = "  +2.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "675.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.12%  "
$ws.Range("E34").Value = "  +15.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "67.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "42.86"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.68%  "
$ws.Range("B37").Value = "ThetaToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +13.50%  "
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.429"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.84%  "
$ws.Range("D39").Value = "0.0₃0855"
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("E40").Value = "  +0.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0496"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.161"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +14.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.44%  "
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.45%  "
